$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.NumberFormat = $fmt
}

# ---------------------------------------------------------------------------
# Step 1: push the footer row (old row 18) down to row 19, keeping its format
# ---------------------------------------------------------------------------
$ws.Range("A18:Q18").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4122)
$ws.Rows.Item(19).RowHeight = 16.5

Set-TextValue $ws.Range("A19") "Monday, 22 September, 2025 10:52 AM"
Set-TextValue $ws.Range("G19") "1/1"
Set-TextValue $ws.Range("K19") "developed by : Abdelaziz Talaat"

$ws.Range("A19:F19").Merge()
$ws.Range("G19:I19").Merge()
$ws.Range("K19:Q19").Merge()

# ---------------------------------------------------------------------------
# Step 2: push the totals row (old row 17) down to row 18, updating the total
# ---------------------------------------------------------------------------
$ws.Range("P17:Q17").Copy()
$ws.Range("P18:Q18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = 24.75
$ws.Range("P18").Value2 = 844
$ws.Range("Q18").Value2 = $null

$ws.Range("P18:Q18").Merge()

# ---------------------------------------------------------------------------
# Step 3: turn row 17 into a new item row (carrying the old row 16 contents:
# "مناديل سولو سحب") using row 16 as the formatting template
# ---------------------------------------------------------------------------
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 25.5

$ws.Range("A17").Value2 = 11
Set-TextValue $ws.Range("C17") "مناديل سولو سحب"
Set-TextValue $ws.Range("H17") "20:0"
Set-TextValue $ws.Range("L17") "0"
Set-TextValue $ws.Range("N17") "45.00"
Set-TextValue $ws.Range("P17") "45.0000"
Set-TextValue $ws.Range("Q17") "1"

$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# ---------------------------------------------------------------------------
# Step 4: update row 16 in place with the new "PLEGICA" item
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("C16") "PLEGICA 1% EYE DROPS 10 ML"
Set-TextValue $ws.Range("H16") "5:0"
Set-TextValue $ws.Range("L16") "1"
Set-TextValue $ws.Range("N16") "27.00"
Set-TextValue $ws.Range("P16") "27.0000"
